# Apply the FHIR IG terminology/metadata corrections described in the commit:
#  - "Experimental" row (B7) on the Metadata sheet gets an explicit text value "false"
#  - "Date" row (B8) on the Metadata sheet is bumped to the new generation timestamp

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Metadata")

# B7 ("Experimental") currently has no value. We need the literal text "false"
# (not the Boolean FALSE that a plain Value assignment of "false" would produce).
# Build it via a formula that concatenates to a string, then convert the formula
# result to a static value in place - this keeps the cell's existing style/format
# and stores it as a normal (shared) text string rather than a boolean.
$cellB7 = $ws.Range("B7")
$cellB7.Formula = "=""fal""&""se"""
$cellB7.Copy() | Out-Null
$cellB7.PasteSpecial(-4163) | Out-Null   # xlPasteValues

# B8 ("Date") simply gets updated to the new ISO-8601 timestamp text.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
